$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.059843338016958
$ws.Cells.Item(2, 4).Value = 1.057488595582969
$ws.Cells.Item(2, 5).Value = 1.064954358630998
$ws.Cells.Item(2, 6).Value = 1.074008256550537
$ws.Cells.Item(2, 9).Value = 1.051316569052303
$ws.Cells.Item(2, 10).Value = 1.064827267394267
$ws.Cells.Item(2, 11).Value = 1.060223427134904
$ws.Cells.Item(2, 12).Value = 1.067668902807051
$ws.Cells.Item(2, 13).Value = 1.076698606907402
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061114870383562
$ws.Cells.Item(3, 4).Value = 1.05845684050984
$ws.Cells.Item(3, 5).Value = 1.06612946171137
$ws.Cells.Item(3, 6).Value = 1.075354904110726
$ws.Cells.Item(3, 9).Value = 1.051749249178972
$ws.Cells.Item(3, 10).Value = 1.065750676083208
$ws.Cells.Item(3, 11).Value = 1.06100512603514
$ws.Cells.Item(3, 12).Value = 1.068658410990063
$ws.Cells.Item(3, 13).Value = 1.077860996936937
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.061936978308239
$ws.Cells.Item(4, 4).Value = 1.059082677099288
$ws.Cells.Item(4, 5).Value = 1.066889439369922
$ws.Cells.Item(4, 6).Value = 1.076226201130027
$ws.Cells.Item(4, 9).Value = 1.052027541051118
$ws.Cells.Item(4, 10).Value = 1.066347004910497
$ws.Cells.Item(4, 11).Value = 1.061509636208021
$ws.Cells.Item(4, 12).Value = 1.069297718326626
$ws.Cells.Item(4, 13).Value = 1.078612520720356
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.062282438033621
$ws.Cells.Item(5, 4).Value = 1.059345617109901
$ws.Cells.Item(5, 5).Value = 1.067208842431637
$ws.Cells.Item(5, 6).Value = 1.076592479655903
$ws.Cells.Item(5, 9).Value = 1.052144133647642
$ws.Cells.Item(5, 10).Value = 1.066597421605677
$ws.Cells.Item(5, 11).Value = 1.061721422496386
$ws.Cells.Item(5, 12).Value = 1.06956625302288
$ws.Cells.Item(5, 13).Value = 1.078928315460299
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.062340433278
$ws.Cells.Item(6, 4).Value = 1.059389756463154
$ws.Cells.Item(6, 5).Value = 1.067262466313841
$ws.Cells.Item(6, 6).Value = 1.076653978654047
$ws.Cells.Item(6, 9).Value = 1.052163686558265
$ws.Cells.Item(6, 10).Value = 1.066639451317128
$ws.Cells.Item(6, 11).Value = 1.06175696423051
$ws.Cells.Item(6, 12).Value = 1.069611327751985
$ws.Cells.Item(6, 13).Value = 1.078981330310139
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.061941594960863
$ws.Cells.Item(7, 4).Value = 1.059086191149768
$ws.Cells.Item(7, 5).Value = 1.06689370760593
$ws.Cells.Item(7, 6).Value = 1.076231095415828
$ws.Cells.Item(7, 9).Value = 1.052029100542529
$ws.Cells.Item(7, 10).Value = 1.066350352087945
$ws.Cells.Item(7, 11).Value = 1.061512467321864
$ws.Cells.Item(7, 12).Value = 1.069301307401733
$ws.Cells.Item(7, 13).Value = 1.078616740954962
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.060273196622993
$ws.Cells.Item(8, 4).Value = 1.057815960874313
$ws.Cells.Item(8, 5).Value = 1.065351572985682
$ws.Cells.Item(8, 6).Value = 1.074463379512886
$ws.Cells.Item(8, 9).Value = 1.051463144105225
$ws.Cells.Item(8, 10).Value = 1.065139582965402
$ws.Cells.Item(8, 11).Value = 1.060487876428537
$ws.Cells.Item(8, 12).Value = 1.068003514067317
$ws.Cells.Item(8, 13).Value = 1.077091572086989
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.057328077904751
$ws.Cells.Item(9, 4).Value = 1.055572346752931
$ws.Cells.Item(9, 5).Value = 1.062631017901786
$ws.Cells.Item(9, 6).Value = 1.071347737967042
$ws.Cells.Item(9, 9).Value = 1.05045292942184
$ws.Cells.Item(9, 10).Value = 1.062996931572487
$ws.Cells.Item(9, 11).Value = 1.058672371640555
$ws.Cells.Item(9, 12).Value = 1.065709100971357
$ws.Cells.Item(9, 13).Value = 1.074399160397423
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.055360974121482
$ws.Cells.Item(10, 4).Value = 1.054072925017169
$ws.Cells.Item(10, 5).Value = 1.060815063024424
$ws.Cells.Item(10, 6).Value = 1.069269983939146
$ws.Cells.Item(10, 9).Value = 1.049770687611687
$ws.Cells.Item(10, 10).Value = 1.061562227159308
$ws.Cells.Item(10, 11).Value = 1.057455170018461
$ws.Cells.Item(10, 12).Value = 1.064174281820858
$ws.Cells.Item(10, 13).Value = 1.072600777576274
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.054508275049999
$ws.Cells.Item(11, 4).Value = 1.053422760687492
$ws.Cells.Item(11, 5).Value = 1.060028165110829
$ws.Cells.Item(11, 6).Value = 1.068370094349853
$ws.Cells.Item(11, 9).Value = 1.049473173482858
$ws.Cells.Item(11, 10).Value = 1.060939465915976
$ws.Cells.Item(11, 11).Value = 1.0569264545258
$ws.Cells.Item(11, 12).Value = 1.063508419673574
$ws.Cells.Item(11, 13).Value = 1.071821203632849
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.054191400859505
$ws.Cells.Item(12, 4).Value = 1.053181122469866
$ws.Cells.Item(12, 5).Value = 1.059735786291483
$ws.Cells.Item(12, 6).Value = 1.068035800214284
$ws.Cells.Item(12, 9).Value = 1.049362346626016
$ws.Cells.Item(12, 10).Value = 1.060707912804527
$ws.Cells.Item(12, 11).Value = 1.056729814501057
$ws.Cells.Item(12, 12).Value = 1.06326089431755
$ws.Cells.Item(12, 13).Value = 1.071531502242107
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.054259378020876
$ws.Cells.Item(13, 4).Value = 1.053232960980351
$ws.Cells.Item(13, 5).Value = 1.059798506659346
$ws.Cells.Item(13, 6).Value = 1.068107509088996
$ws.Cells.Item(13, 9).Value = 1.049386133726242
$ws.Cells.Item(13, 10).Value = 1.060757592269731
$ws.Cells.Item(13, 11).Value = 1.056772005873752
$ws.Cells.Item(13, 12).Value = 1.063313998195261
$ws.Cells.Item(13, 13).Value = 1.071593650242807
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.05448208508726
$ws.Cells.Item(14, 4).Value = 1.053402789630372
$ws.Cells.Item(14, 5).Value = 1.060003998838503
$ws.Cells.Item(14, 6).Value = 1.068342462229674
$ws.Cells.Item(14, 9).Value = 1.049464018973031
$ws.Cells.Item(14, 10).Value = 1.060920330409482
$ws.Cells.Item(14, 11).Value = 1.056910205344776
$ws.Cells.Item(14, 12).Value = 1.063487963135149
$ws.Cells.Item(14, 13).Value = 1.071797259556505
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.054619283086456
$ws.Cells.Item(15, 4).Value = 1.053507408280012
$ws.Cells.Item(15, 5).Value = 1.060130597345736
$ws.Cells.Item(15, 6).Value = 1.068487219853426
$ws.Cells.Item(15, 9).Value = 1.049511964611079
$ws.Cells.Item(15, 10).Value = 1.061020567941765
$ws.Cells.Item(15, 11).Value = 1.056995321198872
$ws.Cells.Item(15, 12).Value = 1.063595122815063
$ws.Cells.Item(15, 13).Value = 1.071922692260882
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.055417545008494
$ws.Cells.Item(16, 4).Value = 1.054116055023059
$ws.Cells.Item(16, 5).Value = 1.060867274413096
$ws.Cells.Item(16, 6).Value = 1.069329701935846
$ws.Cells.Item(16, 9).Value = 1.049790388292821
$ws.Cells.Item(16, 10).Value = 1.061603525410081
$ws.Cells.Item(16, 11).Value = 1.057490223973508
$ws.Cells.Item(16, 12).Value = 1.064218445781366
$ws.Cells.Item(16, 13).Value = 1.072652496876572
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.055918021743466
$ws.Cells.Item(17, 4).Value = 1.054497599179297
$ws.Cells.Item(17, 5).Value = 1.061329215657172
$ws.Cells.Item(17, 6).Value = 1.069858110242154
$ws.Cells.Item(17, 9).Value = 1.049964473097653
$ws.Cells.Item(17, 10).Value = 1.061968789047269
$ws.Cells.Item(17, 11).Value = 1.057800217598442
$ws.Cells.Item(17, 12).Value = 1.064609096322195
$ws.Cells.Item(17, 13).Value = 1.073110050792981
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.056209851884329
$ws.Cells.Item(18, 4).Value = 1.054720060284931
$ws.Cells.Item(18, 5).Value = 1.06159860266589
$ws.Cells.Item(18, 6).Value = 1.070166302017517
$ws.Cells.Item(18, 9).Value = 1.050065811407155
$ws.Cells.Item(18, 10).Value = 1.062181694218566
$ws.Cells.Item(18, 11).Value = 1.057980871705878
$ws.Cells.Item(18, 12).Value = 1.064836833320464
$ws.Cells.Item(18, 13).Value = 1.073376851279344
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.056309343366601
$ws.Cells.Item(19, 4).Value = 1.054795899060566
$ws.Cells.Item(19, 5).Value = 1.061690447354269
$ws.Cells.Item(19, 6).Value = 1.070271384194956
$ws.Cells.Item(19, 9).Value = 1.050100330865073
$ws.Cells.Item(19, 10).Value = 1.062254264541902
$ws.Cells.Item(19, 11).Value = 1.058042443029685
$ws.Cells.Item(19, 12).Value = 1.064914465000923
$ws.Cells.Item(19, 13).Value = 1.073467809378153
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.055864334611333
$ws.Cells.Item(20, 4).Value = 1.054456672133725
$ws.Cells.Item(20, 5).Value = 1.061279659498364
$ws.Cells.Item(20, 6).Value = 1.069801419113466
$ws.Cells.Item(20, 9).Value = 1.049945816393657
$ws.Cells.Item(20, 10).Value = 1.06192961493448
$ws.Cells.Item(20, 11).Value = 1.057766974775885
$ws.Cells.Item(20, 12).Value = 1.064567195945602
$ws.Cells.Item(20, 13).Value = 1.07306096819575
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.054416507408872
$ws.Cells.Item(21, 4).Value = 1.053352783160342
$ws.Cells.Item(21, 5).Value = 1.059943489020394
$ws.Cells.Item(21, 6).Value = 1.068273275384355
$ws.Cells.Item(21, 9).Value = 1.049441092468313
$ws.Cells.Item(21, 10).Value = 1.060872414506236
$ws.Cells.Item(21, 11).Value = 1.056869516007457
$ws.Cells.Item(21, 12).Value = 1.063436740189037
$ws.Cells.Item(21, 13).Value = 1.071737305384715
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.053505367050572
$ws.Cells.Item(22, 4).Value = 1.05265792340598
$ws.Cells.Item(22, 5).Value = 1.059102863924545
$ws.Cells.Item(22, 6).Value = 1.067312265963715
$ws.Cells.Item(22, 9).Value = 1.049121918476128
$ws.Cells.Item(22, 10).Value = 1.060206367797926
$ws.Cells.Item(22, 11).Value = 1.056303791534092
$ws.Cells.Item(22, 12).Value = 1.062724850983608
$ws.Cells.Item(22, 13).Value = 1.07090429534345
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.053988460194385
$ws.Cells.Item(23, 4).Value = 1.053026358398878
$ws.Cells.Item(23, 5).Value = 1.059548545693479
$ws.Cells.Item(23, 6).Value = 1.06782173591449
$ws.Cells.Item(23, 9).Value = 1.049291293002798
$ws.Cells.Item(23, 10).Value = 1.060559579975782
$ws.Cells.Item(23, 11).Value = 1.056603831606477
$ws.Cells.Item(23, 12).Value = 1.063102344700493
$ws.Cells.Item(23, 13).Value = 1.071345963906493
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.055888593791706
$ws.Cells.Item(24, 4).Value = 1.054475165573181
$ws.Cells.Item(24, 5).Value = 1.061302051964834
$ws.Cells.Item(24, 6).Value = 1.069827035455501
$ws.Cells.Item(24, 9).Value = 1.049954247180501
$ws.Cells.Item(24, 10).Value = 1.061947316483723
$ws.Cells.Item(24, 11).Value = 1.057781996270579
$ws.Cells.Item(24, 12).Value = 1.064586129301717
$ws.Cells.Item(24, 13).Value = 1.073083146764972
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.058090097352469
$ws.Cells.Item(25, 4).Value = 1.056153015018199
$ws.Cells.Item(25, 5).Value = 1.063334731977654
$ws.Cells.Item(25, 6).Value = 1.072153308375985
$ws.Cells.Item(25, 9).Value = 1.050715633650705
$ws.Cells.Item(25, 10).Value = 1.06355195336884
$ws.Cells.Item(25, 11).Value = 1.059142924953388
$ws.Cells.Item(25, 12).Value = 1.066303169535266
$ws.Cells.Item(25, 13).Value = 1.075095807636301
